$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tabel")

# Replace formula-driven boolean values in columns B and C (rows 2-30)
# with hardcoded TRUE (1) values, removing the TRUE()/FALSE() formulas.
for ($r = 2; $r -le 30; $r++) {
    $ws.Cells.Item($r, 2).Value = $true
    $ws.Cells.Item($r, 3).Value = $true
}

# Update the active selection on the sheet to C3:C30 with active cell C3.
$ws.Activate()
$ws.Range("C3:C30").Select()
